$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.936.26"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.09"
$ws.Range("E3").Value = "  -2.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.32"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4978"
$ws.Range("E7").Value = "  -3.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  -4.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09043"
$ws.Range("E9").Value = "  -7.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.56"
$ws.Range("E11").Value = "  -1.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.298"
$ws.Range("E12").Value = "  -3.64%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.65"
$ws.Range("E13").Value = "  -2.52%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.878.68"
$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.214"
$ws.Range("E15").Value = "  -3.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("E17").Value = "  -3.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.75"
$ws.Range("E18").Value = "  -4.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06649"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.87"
$ws.Range("E20").Value = "  -2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.122"
$ws.Range("E22").Value = "  -3.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.996.47"
$ws.Range("E23").Value = "  -2.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.290"
$ws.Range("E25").Value = "  -0.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.105.51"
$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.518"
$ws.Range("E27").Value = "  -6.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.88"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.70"
$ws.Range("E29").Value = "  -2.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.55"
$ws.Range("E30").Value = "  -1.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1053"
$ws.Range("E31").Value = "  -2.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.054"
$ws.Range("E32").Value = "  -4.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.566"
$ws.Range("E33").Value = "  -3.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.588"
$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.282"
$ws.Range("E35").Value = "  -6.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06531"
$ws.Range("E36").Value = "  -3.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02393"
$ws.Range("E37").Value = "  -1.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2177"
$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.272"
$ws.Range("E39").Value = "  +6.56%  "

$ws.Range("E40").Value = "  -6.41%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.60"
$ws.Range("E41").Value = "  -1.96%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6365"
$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.895"
$ws.Range("E43").Value = "  -3.96%  "

$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("E45").Value = "  -2.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5987"
$ws.Range("E46").Value = "  -1.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.285"
$ws.Range("E47").Value = "  +0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.676"
$ws.Range("E48").Value = "  -2.84%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.215"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.970"
$ws.Range("E50").Value = "  -3.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.45"
$ws.Range("E51").Value = "  -3.75%  "
